$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Terry Rozier", "PG", "Miami Heat"),
    @("Dalton Knecht", "SG,SF", "Los Angeles Lakers"),
    @("Donte DiVincenzo", "SG,SF", "Minnesota Timberwolves"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Nicolas Claxton", "C", "Brooklyn Nets")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
